$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B..G to C..H)
$ws.Columns("B:B").Insert()

# Insert a new row before row 1 (shifts old rows 1..5 to 2..6)
$ws.Rows("1:1").Insert()

# Row 2 (old header row, now shifted down): rename the degree label cell to ModuleName,
# and add the FinalModuleMark header next to it.
$ws.Range("A2").Value = "ModuleName"
$ws.Range("B2").Value = "FinalModuleMark"

# New row 1: degree label + overall average formula (bold, matching the header row style)
$ws.Range("A1").Value = "Degree: BComHons Information Systems Management"
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Formula = "=SUM(B3:B6)/COUNTA(B3:B6)"
$ws.Range("B1").Font.Bold = $true

# Column B rows 3-6: final module mark formulas (row4:6 share one formula)
$ws.Range("B3").Formula = "=C3*D3+E3*F3+G3*H3"
$ws.Range("B4:B6").Formula = "=C4*D4+E4*F4+G4*H4"

# Data corrections accompanying the new layout
$ws.Range("E4").Value = 0
$ws.Range("C6").Value = 0

# Column widths (bestFit-equivalent width for the new FinalModuleMark header column)
$ws.Range("B:B").ColumnWidth = 14.8

# Selection / active cell
[void]$ws.Range("A12").Select()
